$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels (accent/casing normalization)
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "potencia"
$ws.Range("C1").Value = "tensao"

# Update voltage ("tensao") values for the three circuits
$ws.Range("C2").Value = 220
$ws.Range("C3").Value = 380
$ws.Range("C4").Value = 220

# Move the active selection
$ws.Range("E9").Select()
